$wb = $excel.ActiveWorkbook

# --- Sheet "Tuning" ---
$tuning = $wb.Worksheets.Item("Tuning")
$tuning.Range("H5").Value = 250
$tuning.Range("J5").Value = "x"

# --- Sheet "Values" ---
$values = $wb.Worksheets.Item("Values")
$values.Range("J3:K4").HorizontalAlignment = -4152  # xlRight
$values.Range("I3").Value = 300
$values.Range("I4").Value = 30

$values.Range("J5:K10").HorizontalAlignment = -4131  # xlLeft
$values.Range("I5").Value = "Ghastly 5-8"
$values.Range("I6").Value = "Liz Fox 8-9"
$values.Range("I7").Value = "Gh Fox 7-10"
$values.Range("I8").Value = "Liz 3-8"
$values.Range("I9").Value = "Fox Liz 2-6"
$values.Range("I10").Value = "x"
$values.Range("J10").Value = 2
$values.Range("I11").Value = "Scarecrow on 2"

# --- Sheet "TODO" ---
$todo = $wb.Worksheets.Item("TODO")
$todo.Range("H4").Value = "toggle countdown on/off, alt. Disabled automatically on difficulty lvl 3"
